$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New "Column with int" column (F) - header text plus two real numeric
# cells (previously everything on this sheet was shared-string text;
# this is the "fix int cell issue" from the commit message: ints are now
# written as actual numbers instead of strings).
$ws1.Cells.Item(1, 6).Value = "Column with int"
$ws1.Cells.Item(2, 6).Value = 10
$ws1.Cells.Item(3, 6).Value = 12

# Sheet1 becomes the active sheet/tab, with G5 selected (one cell to the
# right of the new column, mirroring the previous D5 selection that sat
# one cell right of the old last column E).
$ws1.Activate() | Out-Null
$ws1.Range("G5").Select() | Out-Null
